$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns: C (fly_cost) and F (day cost per all people) ---
$ws.Columns("C:C").Insert()
$ws.Columns("F:F").Insert()

# --- Insert eleven new rows after row 2 to make room for the daily detail
#     rows (new rows 3-11), plus the relocated "total" row (12) and the two
#     grand-total rows (13-14) ---
$ws.Rows("3:13").Insert()

# --- Header row (row 1) ---
$ws.Range("C1").Value = "fly_cost"
$ws.Range("F1").Value = "day cost per all people"

# --- Daily detail rows 2-11 ---
# Row 2
$ws.Range("A2").Value = "'2021-06-16"
$ws.Range("B2").Value = 6.1
$ws.Range("C2").Value = 133.1
$ws.Range("D2").Value = 3787
$ws.Range("E2").Value = 38.4
$ws.Range("F2").Value = 7.5
$ws.Range("G2").Value = "transit to Greenland. bring any cargo from Iceland?"

# Row 3
$ws.Range("A3").Value = "'2021-06-17"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 43.7
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 7.5
$ws.Range("G3").Value = "transit to Greenland. bring any cargo from Iceland?"

# Row 4
$ws.Range("A4").Value = "'2021-06-17"
$ws.Range("B4").Value = 2.8
$ws.Range("C4").Value = 61.1
$ws.Range("D4").Value = 2414
$ws.Range("E4").Value = 24.5
$ws.Range("F4").Value = 15
$ws.Range("G4").Value = "return to SFJ to overnight"

# Row 5
$ws.Range("A5").Value = "'2021-06-18"
$ws.Range("B5").Value = 2.3
$ws.Range("C5").Value = 50.5
$ws.Range("D5").Value = 1993
$ws.Range("E5").Value = 20.2
$ws.Range("F5").Value = 15
$ws.Range("G5").Value = "return to SFJ to overnight"

# Row 6
$ws.Range("A6").Value = "'2021-06-19"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 44.6
$ws.Range("D6").Value = 1363
$ws.Range("E6").Value = 13.8
$ws.Range("F6").Value = 15
$ws.Range("G6").Value = "return to JAV to overnight"

# Row 7
$ws.Range("A7").Value = "'2021-06-21"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 43.7
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 15
$ws.Range("G7").Value = "return to JAV to overnight"

# Row 8
$ws.Range("A8").Value = "'2021-06-21"
$ws.Range("B8").Value = 3.7
$ws.Range("C8").Value = 81.7
$ws.Range("D8").Value = 3213
$ws.Range("E8").Value = 32.6
$ws.Range("F8").Value = 15
$ws.Range("G8").Value = "to UAK to overnight"

# Row 9
$ws.Range("A9").Value = "'2021-06-22"
$ws.Range("B9").Value = 1.8
$ws.Range("C9").Value = 39.2
$ws.Range("D9").Value = 718
$ws.Range("E9").Value = 7.3
$ws.Range("F9").Value = 15
$ws.Range("G9").Value = "transit"

# Row 10
$ws.Range("A10").Value = "'2021-06-23"
$ws.Range("B10").Value = 2.9
$ws.Range("C10").Value = 64.09999999999999
$ws.Range("D10").Value = 1940
$ws.Range("E10").Value = 19.7
$ws.Range("F10").Value = 15
$ws.Range("G10").Value = "transit to KUS"

# Row 11
$ws.Range("A11").Value = "'2021-06-24"
$ws.Range("B11").Value = 3.4
$ws.Range("C11").Value = 73.7
$ws.Range("D11").Value = 1350
$ws.Range("E11").Value = 13.7
$ws.Range("F11").Value = 15
$ws.Range("G11").Value = "transit to AEY"

# --- Totals row (12) ---
$ws.Range("A12").Value = "total"
$ws.Range("B12").Value = 29
$ws.Range("C12").Value = 635.4
$ws.Range("D12").Value = 16778
$ws.Range("E12").Value = 170.2
$ws.Range("F12").Value = 135

# --- Grand total rows (13-14) ---
$ws.Range("A13").Value = "grand total (MDKK)"
$ws.Range("B13").Value = 0.9405999999999999

$ws.Range("A14").Value = "grand total incl. quarantine (MDKK)"
$ws.Range("B14").Value = 1.0156
